$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the newly-populated row 17 of the results table
$ws.Range("B17").Value = 13
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = "Sigmoid (10)"
$ws.Range("E17").Value = "ReLU (8)"
$ws.Range("H17").Value = "Sigmoid"
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 0.209
$ws.Range("L17").Value = "Test size reduced from 0.33 to 0.2"

# Update the view state to match the author's saved selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("B18").Select()
